# Script scraped data update for israel_ligat-ha-al_2023-2024.xlsx
# - Swaps the F:V (match detail) content between several pairs of rows that
#   had been re-ordered by the upstream scraper.
# - Appends a new match row (row 98) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowDetails {
    param($rowA, $rowB)
    $rangeA = $ws.Range("F$($rowA):V$($rowA)")
    $rangeB = $ws.Range("F$($rowB):V$($rowB)")
    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()
    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Pairs of rows whose F:V (home..url) details were swapped back into the
# correct chronological order.
Swap-RowDetails 4 5
Swap-RowDetails 31 32
Swap-RowDetails 44 45
Swap-RowDetails 58 60
Swap-RowDetails 72 73

# Append the new match row (row 98) using row 97 as the formatting template.
$ws.Range("A97:V97").Copy()
$ws.Range("A98:V98").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A98").Value = 97
$ws.Range("B98").Value = "israel"
$ws.Range("C98").Value = "ligat-ha-al"
$ws.Range("D98").Value = "2023-2024"
$ws.Range("E98").Value = 45295.8125
$ws.Range("F98").Value = "Hapoel Haifa"
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = "Beitar Jerusalem"
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = 2.52
$ws.Range("K98").Value = "31/12/2024 19:17"
$ws.Range("L98").Value = 2.51
$ws.Range("M98").Value = "04/01/2024 19:26"
$ws.Range("N98").Value = 3.23
$ws.Range("O98").Value = "31/12/2024 19:17"
$ws.Range("P98").Value = 3.43
$ws.Range("Q98").Value = "04/01/2024 19:26"
$ws.Range("R98").Value = 2.74
$ws.Range("S98").Value = "31/12/2024 19:17"
$ws.Range("T98").Value = 2.82
$ws.Range("U98").Value = "04/01/2024 19:26"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-haifa-beitar-jerusalem/OtKnISc3/"
